# Remove three example rows (old rows 3, 5 and 6) from the firewall rule
# list sheet. Excel shifts the remaining rows up, automatically renumbering
# the sheetData rows, shrinking the dimension, adjusting the data
# validation ranges, and dropping any shared strings that become unused
# (the "0.0.0.0/0", "192.168.254.0/27" and "80, 443" entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 3 first ...
$ws.Rows(3).Delete()
# ... then the two rows that are now rows 4 and 5 (previously rows 5 and 6).
$ws.Rows(4).Delete()
$ws.Rows(4).Delete()

# Leave the selection on the row that now holds the "Outbound /
# 192.168.12.0/24, 192.168.14.0/27 / 9000" rule, matching the author's
# final selection state (whole row 4 selected).
[void]$ws.Range("A4:XFD4").Select()
